$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EditViewTest")

$ws.Range("A2").Value = "testView"
$ws.Range("B2").Value = "testViewEdited"

$ws.Activate()
$ws.Range("B2").Select()
